# Insert a new data row at row 22 (shifts existing rows 22:101 down to 23:102)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(22).Insert()

$ws.Range("A22").Value = 10
$ws.Range("B22").Value = "Vega Modelo de Temuco"
$ws.Range("C22").Value = "La Araucanía"
$ws.Range("D22").Value = 45243
$ws.Range("E22").Value = 9
$ws.Range("F22").Value = 100112026
$ws.Range("G22").Value = "Haba"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 65
$ws.Range("K22").Value = 10000
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = 10000
$ws.Range("N22").Value = "$/saco 25 kilos"
$ws.Range("O22").Value = "Región del Maule"
$ws.Range("P22").Value = 400
$ws.Range("Q22").Value = 25
$ws.Range("R22").Value = "Hortaliza"
